$d = $word.ActiveDocument

$replacements = @(
    @("572×4=2288", "295×2=590"),
    @("537×5=2685", "113×6=678"),
    @("416×6=2496", "834×3=2502"),
    @("726×5=3630", "647×3=1941"),
    @("262×7=1834", "157×6=942"),
    @("867×7=6069", "634×7=4438"),
    @("378×3=1134", "841×7=5887"),
    @("447×2=894", "814×4=3256"),
    @("281×4=1124", "209×4=836"),
    @("574×4=2296", "481×7=3367"),
    @("858×7=6006", "439×5=2195"),
    @("808×8=6464", "359×9=3231"),
    @("624×4=2496", "226×9=2034"),
    @("637×7=4459", "983×4=3932"),
    @("945×4=3780", "586×7=4102"),
    @("766×6=4596", "683×2=1366"),
    @("367×7=2569", "583×8=4664"),
    @("706×8=5648", "440×6=2640"),
    @("776×3=2328", "897×7=6279"),
    @("550×5=2750", "683×2=1366"),
    @("358×8=2864", "516×8=4128"),
    @("702×6=4212", "677×6=4062"),
    @("182×4=728", "649×3=1947"),
    @("754×5=3770", "798×3=2394"),
    @("644×9=5796", "913×3=2739")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
